$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: average of the |S*|/n column (J) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# --- Rows 14-17: summary labels + aggregate formulas ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Formatting for the summary values: bold, size 12, vertically centered ---
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

# Clone the formatting to the rest of the summary cells with a format-only
# paste, so no stray/duplicate cell styles get created in the stylesheet.
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection mirrors what was highlighted when the file was saved ---
$ws.Range("A14:B17").Select()

# --- Print setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
